$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.068.08"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "2.763.19"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.71"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.32"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -3.22%  "
$ws.Range("E9").Value = "  -3.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.81"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -14.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.386"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("D13").Value = "3.248.74"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.94"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").Value = "63.678.55"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D17").Value = "2.764.87"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.12"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.84"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -3.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "356.09"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -3.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -5.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.527"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -6.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.05"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -3.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.170"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -3.81%  "
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  -6.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.32"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.96"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -3.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.27"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.95"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("E33").Value = "  -3.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.15"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.80"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "350.34"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.28"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("E41").Value = "  -2.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.10"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.45"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -4.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.82"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -4.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0587"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -4.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "137.58"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.630"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("E48").Value = "  -2.95%  "
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.03"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -0.11%  "
